$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text representation (e.g. leading/trailing zeros,
# percent-style strings, multi-dot numbers) instead of being auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.860.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.650.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.69'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3891'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3822'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.86'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.340'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08450'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.90'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.013'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.022'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001312'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.652.62'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.06'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06983'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.56'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.962'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.66'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.862.79'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.447'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.936'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.02'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.22'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.415'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.43'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.752'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.483'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.832.70'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08126'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9986'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.656'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02903'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2675'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.69'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09110'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7554'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.52'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.42'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6936'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.440'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.104'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08280'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.60'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.228'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.87%  '
